# NC92Soil - Versione 0.6 bis
# - Implementata gestione ottimizzata delle permutazioni in modo da calcolarle
#   una sola volta per ogni cluster
# - Corretto bug che provocava la creazione di piu' export identici nel caso di
#   DH se il foglio Soils conteneva piu' percentili per la Vs

$wb = $excel.ActiveWorkbook

$wsClusters = $wb.Worksheets.Item("Clusters")
$wsProfiles = $wb.Worksheets.Item("Profiles")

# ---------------------------------------------------------------------------
# Clusters sheet: add a 3rd sub-cluster (SUB3) row for the IW1 cluster, and
# fix the brick-thickness of SUB2 (5 -> 4).
# ---------------------------------------------------------------------------

# Bug fix: SUB2 brick thickness corrected from 5 to 4
$wsClusters.Range("D3").Value = 4

# ---------------------------------------------------------------------------
# Profiles sheet: the export-permutation columns are renamed (P1/P2/P3 ->
# DH1/SOND1/SOND2) and a 4th column (DH2) is introduced because the Soils
# sheet can now carry more than one Vs percentile per cluster.
# ---------------------------------------------------------------------------

# Column C (was P2) keeps the same kind of data but gets a new header and a
# couple of its rows change cluster letter (B -> G).
$wsProfiles.Range("C3").Value = "G;3"
$wsProfiles.Range("A4").Value = "G;8;300"
$wsProfiles.Range("C5").Value = "G;9"

# New row in Clusters describing the extra permutation cluster (SUB3)
$wsClusters.Range("A4").Value = "IW1"
$wsClusters.Range("B4").Value = "SUB3"
$wsClusters.Range("C4").Value = 16
$wsClusters.Range("D4").Value = 3
$wsClusters.Range("E4").Value = "Spettro UHS 2.txt; Spettro UHS 3.txt"
$wsClusters.Range("F4").Value = 33
$wsClusters.Range("G4").Value = 33
$wsClusters.Range("H4").Value = 33

# Re-style the data rows: they used to be horizontally centered only: now
# they are simply vertically centered (left aligned horizontally).
$rngData = $wsClusters.Range("A2:H4")
$rngData.VerticalAlignment = -4108
$rngData.HorizontalAlignment = 1

# New header + data for DH1 (column A keeps header, was "P1")
$wsProfiles.Range("A1").Value = "DH1"

# New column B (DH2) - second percentile permutation for the same cluster
$wsProfiles.Range("B3").Value = "A;5;250"
$wsProfiles.Range("B4").Value = "G;9;320"
$wsProfiles.Range("B5").Value = "A;8;380"

# New headers for the (former P2/P3, now SOND1/SOND2) columns
$wsProfiles.Range("C1").Value = "SOND1"
$wsProfiles.Range("D1").Value = "SOND2"

# New column D (SOND2) data
$wsProfiles.Range("D2").Value = "Spettro UHS 3.txt"
$wsProfiles.Range("D3").Value = "S;5"
$wsProfiles.Range("D4").Value = "A;7"
$wsProfiles.Range("D5").Value = "G;10"

# Header for column B (DH2) - written last on purpose to control the order
# new strings are introduced in.
$wsProfiles.Range("B1").Value = "DH2"

# Match the D1 header formatting to the rest of the header row (bold, centred)
$wsProfiles.Range("D1").Font.Bold = $true
$wsProfiles.Range("D1").HorizontalAlignment = -4108
$wsProfiles.Range("D1").VerticalAlignment = -4108

# Columns B and C now share the same (wider) column width, column D takes on
# the narrower width that column C used to have.
$wsProfiles.Columns.Item(3).ColumnWidth = $wsProfiles.Columns.Item(2).ColumnWidth()
$wsProfiles.Columns.Item(4).ColumnWidth = 14.33

# ---------------------------------------------------------------------------
# Selection / active sheet bookkeeping to match the saved workbook state.
# ---------------------------------------------------------------------------
$wsClusters.Activate() | Out-Null
$wsClusters.Range("H10").Select() | Out-Null

$wsProfiles.Activate() | Out-Null
$wsProfiles.Range("C8").Select() | Out-Null
